# Update countries & provincias Spain
#
# This reproduces (in terms of resolved cell content) a refresh of the
# "Pais" COVID dashboard: a few country rows are relabeled because two
# entries ("Angola" and "Santa Lucia") were re-ranked earlier in the
# underlying source table, the "last updated" timestamp moves from
# 20:20 to 21:37, and the daily totals for a batch of countries are
# refreshed with newer figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 1 de Octubre de 2020 a las 21:37"

# --- Country column re-labelling caused by Angola moving up the ranking,
#     right after Nicaragua (rows 125-128 shift one country down) ---
$ws.Range("A125").Value = "Angola"
$ws.Range("A126").Value = "Hong Kong"
$ws.Range("A127").Value = "Congo"
$ws.Range("A128").Value = "Guinea Ecuatorial"

# --- Santa Lucia moving ahead of Nueva Caledonia (rows 207-208 swap) ---
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("A208").Value = "Nueva Caledonia"

# --- Refreshed case figures ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7480283
$ws.Range("C4").Value = 33001
$ws.Range("D4").Value = 4722657
$ws.Range("E4").Value = 2545301
$ws.Range("G4").Value = 585
$ws.Range("H4").Value = 212325

# Row 5 - India
$ws.Range("B5").Value = 6391960
$ws.Range("C5").Value = 81693
$ws.Range("D5").Value = 5348653
$ws.Range("E5").Value = 943503
$ws.Range("G5").Value = 1096
$ws.Range("H5").Value = 99804

# Row 21 - Turquia
$ws.Range("B21").Value = 320070
$ws.Range("C21").Value = 1407
$ws.Range("D21").Value = 281151
$ws.Range("E21").Value = 30657
$ws.Range("G21").Value = 67
$ws.Range("H21").Value = 8262

# Row 29 - Canada
$ws.Range("B29").Value = 160265
$ws.Range("C29").Value = 1507
$ws.Range("D29").Value = 136080
$ws.Range("E29").Value = 14869

# Row 30 - Ecuador
$ws.Range("B30").Value = 138584
$ws.Range("C30").Value = 1537
$ws.Range("E30").Value = 14855
$ws.Range("G30").Value = 78
$ws.Range("H30").Value = 11433

# Row 60 - Uzbekistan
$ws.Range("B60").Value = 57190
$ws.Range("C60").Value = 473
$ws.Range("D60").Value = 53835
$ws.Range("E60").Value = 2884

# Row 68 - Libano
$ws.Range("D68").Value = 18103
$ws.Range("E68").Value = 22405

# Row 104 - Maldivas
$ws.Range("B104").Value = 10354
$ws.Range("C104").Value = 63
$ws.Range("D104").Value = 9187
$ws.Range("E104").Value = 1133

# Row 113 - Zimbabue
$ws.Range("B113").Value = 7850
$ws.Range("C113").Value = 12
$ws.Range("D113").Value = 6312
$ws.Range("E113").Value = 1310

# Row 117 - Cabo Verde
$ws.Range("B117").Value = 6126
$ws.Range("C117").Value = 102
$ws.Range("D117").Value = 5338
$ws.Range("E117").Value = 727
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 61

# Row 123 - Republica de Yibuti
$ws.Range("B123").Value = 5417
$ws.Range("C123").Value = 1
$ws.Range("D123").Value = 5346
$ws.Range("E123").Value = 10

# Row 125 - Angola (new figures)
$ws.Range("B125").Value = 5114
$ws.Range("C125").Value = 142
$ws.Range("D125").Value = 2082
$ws.Range("E125").Value = 2847
$ws.Range("G125").Value = 2
$ws.Range("H125").Value = 185

# Row 126 - Hong Kong (new figures)
$ws.Range("B126").Value = 5098
$ws.Range("C126").Value = 10
$ws.Range("D126").Value = 4836
$ws.Range("E126").Value = 157
$ws.Range("H126").Value = 105

# Row 127 - Congo (new figures)
$ws.Range("B127").Value = 5089
$ws.Range("D127").Value = 3887
$ws.Range("E127").Value = 1113
$ws.Range("H127").Value = 89

# Row 128 - Guinea Ecuatorial (new figures)
$ws.Range("B128").Value = 5045
$ws.Range("C128").Value = 15
$ws.Range("D128").Value = 4879
$ws.Range("E128").Value = 83
$ws.Range("H128").Value = 83

# Row 130 - Ruanda
$ws.Range("B130").Value = 4843
$ws.Range("C130").Value = 3
$ws.Range("D130").Value = 3181
$ws.Range("E130").Value = 1633

# Row 139 - Somalia
$ws.Range("B139").Value = 3593
$ws.Range("C139").Value = 5
$ws.Range("D139").Value = 3001
$ws.Range("E139").Value = 493

# Row 166 - Republica del Chad
$ws.Range("B166").Value = 1203
$ws.Range("C166").Value = 3
$ws.Range("D166").Value = 1067
$ws.Range("E166").Value = 51

# Row 178 - Comoras
$ws.Range("B178").Value = 484
$ws.Range("C178").Value = 5
$ws.Range("D178").Value = 466
$ws.Range("E178").Value = 11

# Row 205 - Dominica
$ws.Range("B205").Value = 31
$ws.Range("C205").Value = 1
$ws.Range("E205").Value = 7
